$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C7").Value = "Consequuntur fugiat "

$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "07-12-2023"
$ws.Range("C14").Value = "Aut et officia repel"
$ws.Range("D14").Value = "Hic ipsam ab volupta"
$ws.Range("E14").Value = "Est aut deserunt qu"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "07-12-2023"

$ws.Range("E15").Select()
